$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to be treated as text so numeric-looking strings
# (e.g. "1.00", "573.50") are not silently coerced into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# --- Update Price (D) and Volume(1h) (E) for rows 2-26 (no row shift in this range) ---
$ws.Range("D2").Value = '63.001.69'
$ws.Range("E2").Value = '  +5.40%  '
$ws.Range("D3").Value = '2.458.78'
$ws.Range("E3").Value = '  +3.74%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '573.50'
$ws.Range("E5").Value = '  +2.57%  '
$ws.Range("D6").Value = '145.97'
$ws.Range("E6").Value = '  +5.99%  '
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("D8").Value = '0.541'
$ws.Range("E8").Value = '  +2.15%  '
$ws.Range("D9").Value = '2.457.85'
$ws.Range("E9").Value = '  +3.86%  '
$ws.Range("D10").Value = '0.111'
$ws.Range("E10").Value = '  +5.81%  '
$ws.Range("D11").Value = '0.160'
$ws.Range("E11").Value = '  +0.89%  '
$ws.Range("D12").Value = '5.23'
$ws.Range("E12").Value = '  +2.84%  '
$ws.Range("D13").Value = '0.353'
$ws.Range("E13").Value = '  +4.65%  '
$ws.Range("D14").Value = '27.36'
$ws.Range("E14").Value = '  +6.95%  '
$ws.Range("D15").Value = '0.0000177'
$ws.Range("E15").Value = '  +7.83%  '
$ws.Range("D16").Value = '2.858.16'
$ws.Range("E16").Value = '  +2.28%  '
$ws.Range("D17").Value = '62.762.99'
$ws.Range("E17").Value = '  +5.07%  '
$ws.Range("D18").Value = '2.442.84'
$ws.Range("E18").Value = '  +3.41%  '
$ws.Range("D19").Value = '7.90'
$ws.Range("E19").Value = '  -1.22%  '
$ws.Range("D20").Value = '10.99'
$ws.Range("E20").Value = '  +4.98%  '
$ws.Range("D21").Value = '328.41'
$ws.Range("E21").Value = '  +2.21%  '
$ws.Range("D22").Value = '4.13'
$ws.Range("E22").Value = '  +2.27%  '
$ws.Range("D23").Value = '2.04'
$ws.Range("E23").Value = '  +13.03%  '
$ws.Range("D24").Value = '0.998'
$ws.Range("E24").Value = '  -0.42%  '
$ws.Range("D25").Value = '65.67'
$ws.Range("E25").Value = '  +2.32%  '
$ws.Range("D26").Value = '636.24'
$ws.Range("E26").Value = '  +13.86%  '

# --- Insert a new row at 27 for "Binance-PegBSC-USD", shifting rows 27-51 down to 28-52 ---
$ws.Rows.Item(27).Insert()

# --- Remove the now-duplicated last row (the old row 51 "Stellar", now at 52) ---
$ws.Rows.Item(52).Delete()

# --- Populate the new row 27 (Coin/Link/Price/Volume) ---
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = '1.10'
$ws.Range("E27").Value = '  +10.26%  '

# --- Update Price (D) and Volume(1h) (E) for rows 28-51 (B/C already carried down by the row insert/shift) ---
$ws.Range("D28").Value = '8.57'
$ws.Range("E28").Value = '  +5.07%  '
$ws.Range("D29").Value = '0.0₃0986'
$ws.Range("E29").Value = '  +6.96%  '
$ws.Range("D30").Value = '2.537.85'
$ws.Range("E30").Value = '  +2.21%  '
$ws.Range("D31").Value = '8.19'
$ws.Range("E31").Value = '  +2.43%  '
$ws.Range("D32").Value = '1.41'
$ws.Range("E32").Value = '  +8.73%  '
$ws.Range("D33").Value = '0.138'
$ws.Range("E33").Value = '  +6.12%  '
$ws.Range("D34").Value = '1.85'
$ws.Range("E34").Value = '  +4.04%  '
$ws.Range("D35").Value = '1.49'
$ws.Range("E35").Value = '  +5.02%  '
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").Value = '4.76'
$ws.Range("E37").Value = '  +5.09%  '
$ws.Range("D38").Value = '0.374'
$ws.Range("E38").Value = '  +2.06%  '
$ws.Range("D39").Value = '152.94'
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("D40").Value = '5.40'
$ws.Range("E40").Value = '  +8.50%  '
$ws.Range("D41").Value = '18.68'
$ws.Range("E41").Value = '  +2.93%  '
$ws.Range("D42").Value = '2.74'
$ws.Range("E42").Value = '  +14.07%  '
$ws.Range("D43").Value = '1.77'
$ws.Range("E43").Value = '  +8.50%  '
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").Value = '0.0₆0288'
$ws.Range("E45").Value = '  -3.23%  '
$ws.Range("D46").Value = '144.90'
$ws.Range("E46").Value = '  +4.52%  '
$ws.Range("D47").Value = '3.60'
$ws.Range("E47").Value = '  +2.38%  '
$ws.Range("D48").Value = '20.44'
$ws.Range("E48").Value = '  +7.57%  '
$ws.Range("D49").Value = '0.602'
$ws.Range("E49").Value = '  +3.04%  '
$ws.Range("D50").Value = '0.0516'
$ws.Range("E50").Value = '  +3.34%  '
$ws.Range("D51").Value = '12.68'
$ws.Range("E51").Value = '  +8.55%  '

# --- Column A is a static row index (0,1,2,...) independent of the coin-data shift;
#     the row Insert/Delete above shifted it along with everything else, so restore it. ---
$ws.Range("A27").Value = 25
$ws.Range("A28").Value = 26
$ws.Range("A29").Value = 27
$ws.Range("A30").Value = 28
$ws.Range("A31").Value = 29
$ws.Range("A32").Value = 30
$ws.Range("A33").Value = 31
$ws.Range("A34").Value = 32
$ws.Range("A35").Value = 33
$ws.Range("A36").Value = 34
$ws.Range("A37").Value = 35
$ws.Range("A38").Value = 36
$ws.Range("A39").Value = 37
$ws.Range("A40").Value = 38
$ws.Range("A41").Value = 39
$ws.Range("A42").Value = 40
$ws.Range("A43").Value = 41
$ws.Range("A44").Value = 42
$ws.Range("A45").Value = 43
$ws.Range("A46").Value = 44
$ws.Range("A47").Value = 45
$ws.Range("A48").Value = 46
$ws.Range("A49").Value = 47
$ws.Range("A50").Value = 48
$ws.Range("A51").Value = 49
